$wb = $excel.ActiveWorkbook

# --- Carrier sheet: add "onwind" row (new carrier, reuses the #66039c color) ---
$wsCarrier = $wb.Worksheets.Item("Carrier")
$wsCarrier.Range("A12").Value = "onwind"
$wsCarrier.Range("B12").Value = "#66039c"

# --- Generator sheet: rename "diesel 1" -> "wind", retag carriers, update p_nom_max ---
$wsGen = $wb.Worksheets.Item("Generator")
$wsGen.Range("B2").Value = "gas CHP"
$wsGen.Range("G2").Value = 20
$wsGen.Range("A3").Value = "wind"
$wsGen.Range("B3").Value = "onwind"
$wsGen.Range("G3").Value = 120

# --- StorageUnit sheet: add a "battery" row ---
$wsStorage = $wb.Worksheets.Item("StorageUnit")
$wsStorage.Range("A2").Value = "battery"
$wsStorage.Range("B2").Value = "bus 1"
$wsStorage.Range("C2").Value = "AC"
$wsStorage.Range("D2").Value = "True"
$wsStorage.Range("E2").Value = "True"
$wsStorage.Range("F2").Value = 0.95
$wsStorage.Range("G2").Value = 0.95
$wsStorage.Range("I2").Value = 6

# --- selections / active sheet ---
$wsCarrier.Range("A6").Select() | Out-Null
$wsGen.Range("B4").Select() | Out-Null
$wsStorage.Range("B5").Select() | Out-Null
$wsStorage.Activate() | Out-Null
